$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "56.679.67"
$ws.Cells.Item(2, 5).Value = "  -2.25%  "
$ws.Cells.Item(3, 4).Value = "2.991.67"
$ws.Cells.Item(3, 5).Value = "  -4.28%  "
$ws.Cells.Item(4, 4).Value = "0.999"
$ws.Cells.Item(4, 5).Value = "  -0.11%  "
$ws.Cells.Item(5, 4).Value = "497.45"
$ws.Cells.Item(5, 5).Value = "  -4.69%  "
$ws.Cells.Item(6, 4).Value = "135.38"
$ws.Cells.Item(6, 5).Value = "  +0.32%  "
$ws.Cells.Item(7, 4).Value = "0.999"
$ws.Cells.Item(7, 5).Value = "  -0.17%  "
$ws.Cells.Item(8, 4).Value = "2.990.31"
$ws.Cells.Item(8, 5).Value = "  -4.32%  "
$ws.Cells.Item(9, 4).Value = "0.426"
$ws.Cells.Item(9, 5).Value = "  -3.86%  "
$ws.Cells.Item(10, 4).Value = "7.29"
$ws.Cells.Item(10, 5).Value = "  +0.22%  "
$ws.Cells.Item(11, 4).Value = "0.105"
$ws.Cells.Item(11, 5).Value = "  -3.70%  "
$ws.Cells.Item(12, 4).Value = "0.354"
$ws.Cells.Item(12, 5).Value = "  -7.32%  "
$ws.Cells.Item(13, 5).Value = "  +0.08%  "
$ws.Cells.Item(14, 4).Value = "3.501.40"
$ws.Cells.Item(14, 5).Value = "  -4.39%  "
$ws.Cells.Item(15, 4).Value = "25.07"
$ws.Cells.Item(15, 5).Value = "  -1.53%  "
$ws.Cells.Item(16, 4).Value = "56.507.18"
$ws.Cells.Item(16, 5).Value = "  -2.38%  "
$ws.Cells.Item(17, 4).Value = "2.984.63"
$ws.Cells.Item(17, 5).Value = "  -4.47%  "
$ws.Cells.Item(18, 4).Value = "0.0000146"
$ws.Cells.Item(18, 5).Value = "  -4.26%  "
$ws.Cells.Item(19, 4).Value = "'5.90"
$ws.Cells.Item(19, 5).Value = "  +0.94%  "
$ws.Cells.Item(20, 4).Value = "12.38"
$ws.Cells.Item(20, 5).Value = "  -5.11%  "
$ws.Cells.Item(21, 4).Value = "7.77"
$ws.Cells.Item(21, 5).Value = "  -2.70%  "
$ws.Cells.Item(22, 4).Value = "325.03"
$ws.Cells.Item(22, 5).Value = "  -5.65%  "
$ws.Cells.Item(23, 5).Value = "  -0.02%  "
$ws.Cells.Item(24, 4).Value = "0.465"
$ws.Cells.Item(24, 5).Value = "  -8.23%  "
$ws.Cells.Item(25, 4).Value = "61.56"
$ws.Cells.Item(25, 5).Value = "  -10.63%  "
$ws.Cells.Item(26, 4).Value = "0.998"
$ws.Cells.Item(26, 5).Value = "  -0.03%  "
$ws.Cells.Item(27, 4).Value = "0.164"
$ws.Cells.Item(27, 5).Value = "  -2.18%  "
$ws.Cells.Item(28, 4).Value = "0.0₃0890"
$ws.Cells.Item(28, 5).Value = "  -7.15%  "
$ws.Cells.Item(29, 5).Value = "  +0.04%  "
$ws.Cells.Item(30, 4).Value = "6.65"
$ws.Cells.Item(30, 5).Value = "  -1.49%  "
$ws.Cells.Item(31, 4).Value = "6.83"
$ws.Cells.Item(31, 5).Value = "  -0.21%  "
$ws.Cells.Item(32, 4).Value = "1.19"
$ws.Cells.Item(32, 5).Value = "  -3.10%  "
$ws.Cells.Item(33, 4).Value = "1.73"
$ws.Cells.Item(33, 5).Value = "  -6.98%  "
$ws.Cells.Item(34, 4).Value = "19.98"
$ws.Cells.Item(34, 5).Value = "  -7.40%  "
$ws.Cells.Item(35, 4).Value = "'154.30"
$ws.Cells.Item(35, 5).Value = "  -2.30%  "
$ws.Cells.Item(36, 4).Value = "4.52"
$ws.Cells.Item(36, 5).Value = "  -5.40%  "
$ws.Cells.Item(37, 2).Value = "ImmutableX"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(37, 4).Value = "1.29"
$ws.Cells.Item(37, 5).Value = "  -4.93%  "
$ws.Cells.Item(38, 2).Value = "Aptos"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(38, 4).Value = "5.66"
$ws.Cells.Item(38, 5).Value = "  -8.65%  "
$ws.Cells.Item(39, 4).Value = "0.0671"
$ws.Cells.Item(39, 5).Value = "  -2.99%  "
$ws.Cells.Item(40, 4).Value = "23.51"
$ws.Cells.Item(40, 5).Value = "  -4.96%  "
$ws.Cells.Item(41, 4).Value = "3.019.90"
$ws.Cells.Item(41, 5).Value = "  -4.35%  "
$ws.Cells.Item(42, 4).Value = "37.37"
$ws.Cells.Item(42, 5).Value = "  -7.44%  "
$ws.Cells.Item(43, 4).Value = "0.999"
$ws.Cells.Item(43, 5).Value = "  -0.10%  "
$ws.Cells.Item(44, 4).Value = "1.02"
$ws.Cells.Item(44, 5).Value = "  -5.57%  "
$ws.Cells.Item(45, 2).Value = "Stacks"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(45, 4).Value = "1.42"
$ws.Cells.Item(45, 5).Value = "  -1.33%  "
$ws.Cells.Item(46, 2).Value = "Mantle"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(46, 4).Value = "'0.640"
$ws.Cells.Item(46, 5).Value = "  -8.18%  "
$ws.Cells.Item(47, 4).Value = "2.195.13"
$ws.Cells.Item(47, 5).Value = "  -2.85%  "
$ws.Cells.Item(48, 4).Value = "3.58"
$ws.Cells.Item(48, 5).Value = "  -8.41%  "
$ws.Cells.Item(49, 4).Value = "1.94"
$ws.Cells.Item(49, 5).Value = "  +6.64%  "
$ws.Cells.Item(50, 4).Value = "0.0238"
$ws.Cells.Item(50, 5).Value = "  +1.72%  "
$ws.Cells.Item(51, 4).Value = "19.41"
$ws.Cells.Item(51, 5).Value = "  -5.28%  "
